# Fruta / hortaliza, semanal
# Insert 4 new weekly observation rows into the "Naranja" subconjunto sheet,
# right above the existing block (old row 49), pushing the remainder of the
# table down by 4 rows (old rows 49-124 become new rows 53-128).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows before the current row 49 - this shifts rows 49:124
# down to 53:128 and extends the used range / dimension automatically.
$ws.Rows("49:52").Insert()

# --- New row 49 : Fukumoto / Primera --------------------------------------
$ws.Range("A49").Value = 7
$ws.Range("B49").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C49").Value = "Ñuble"
$ws.Range("D49").Value = 44413
$ws.Range("E49").Value = 16
$ws.Range("F49").Value = "Fruta"
$ws.Range("G49").Value = 100102
$ws.Range("H49").Value = "Cítricos"
$ws.Range("I49").Value = 100102005
$ws.Range("J49").Value = "Naranja"
$ws.Range("K49").Value = "Fukumoto"
$ws.Range("L49").Value = "Primera"
$ws.Range("M49").Value = 120
$ws.Range("N49").Value = 6000
$ws.Range("O49").Value = 6200
$ws.Range("P49").Value = 6100
$ws.Range("Q49").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R49").Value = "Región de O'Higgins"
$ws.Range("S49").Value = 407
$ws.Range("T49").Value = 15

# --- New row 50 : Fukumoto / Segunda ---------------------------------------
$ws.Range("A50").Value = 7
$ws.Range("B50").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C50").Value = "Ñuble"
$ws.Range("D50").Value = 44413
$ws.Range("E50").Value = 16
$ws.Range("F50").Value = "Fruta"
$ws.Range("G50").Value = 100102
$ws.Range("H50").Value = "Cítricos"
$ws.Range("I50").Value = 100102005
$ws.Range("J50").Value = "Naranja"
$ws.Range("K50").Value = "Fukumoto"
$ws.Range("L50").Value = "Segunda"
$ws.Range("M50").Value = 100
$ws.Range("N50").Value = 5500
$ws.Range("O50").Value = 5800
$ws.Range("P50").Value = 5650
$ws.Range("Q50").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R50").Value = "Región de O'Higgins"
$ws.Range("S50").Value = 377
$ws.Range("T50").Value = 15

# --- New row 51 : Navel Late / Primera -------------------------------------
$ws.Range("A51").Value = 7
$ws.Range("B51").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C51").Value = "Ñuble"
$ws.Range("D51").Value = 44413
$ws.Range("E51").Value = 16
$ws.Range("F51").Value = "Fruta"
$ws.Range("G51").Value = 100102
$ws.Range("H51").Value = "Cítricos"
$ws.Range("I51").Value = 100102005
$ws.Range("J51").Value = "Naranja"
$ws.Range("K51").Value = "Navel Late"
$ws.Range("L51").Value = "Primera"
$ws.Range("M51").Value = 240
$ws.Range("N51").Value = 6500
$ws.Range("O51").Value = 7000
$ws.Range("P51").Value = 6750
$ws.Range("Q51").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R51").Value = "Región de O'Higgins"
$ws.Range("S51").Value = 450
$ws.Range("T51").Value = 15

# --- New row 52 : Navel Late / Segunda -------------------------------------
$ws.Range("A52").Value = 7
$ws.Range("B52").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C52").Value = "Ñuble"
$ws.Range("D52").Value = 44413
$ws.Range("E52").Value = 16
$ws.Range("F52").Value = "Fruta"
$ws.Range("G52").Value = 100102
$ws.Range("H52").Value = "Cítricos"
$ws.Range("I52").Value = 100102005
$ws.Range("J52").Value = "Naranja"
$ws.Range("K52").Value = "Navel Late"
$ws.Range("L52").Value = "Segunda"
$ws.Range("M52").Value = 160
$ws.Range("N52").Value = 5500
$ws.Range("O52").Value = 6000
$ws.Range("P52").Value = 5750
$ws.Range("Q52").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R52").Value = "Región de O'Higgins"
$ws.Range("S52").Value = 383
$ws.Range("T52").Value = 15
